# Incremental Cost to Produce Substitute Fuel for LCFS — add three new
# reference-fuel rows (heavy/residual fuel oil, LPG propane/butane,
# hydrogen) to the ICtPSFfL sheet, plus a new header label in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ICtPSFfL")

# New rows 9-11: label in column A, "=0" placeholder formulas across B:AK
# (columns 2 through 37), matching the existing rows 4/5/8 pattern.

$ws.Range("A9").Value = "heavy or residual fuel oil"
for ($col = 2; $col -le 37; $col++) {
    $ws.Cells.Item(9, $col).Formula = "=0"
}

$ws.Range("A10").Value = "LPG propane or butane"
for ($col = 2; $col -le 37; $col++) {
    $ws.Cells.Item(10, $col).Formula = "=0"
}

$ws.Range("A11").Value = "hydrogen"
for ($col = 2; $col -le 37; $col++) {
    $ws.Cells.Item(11, $col).Formula = "=0"
}

# New bold header label in A1 ("Cost ($/BTU)"), added last so it lands at
# the end of the shared-string table.
$ws.Range("A1").Value = "Cost (`$/BTU)"
$ws.Range("A1").Font.Bold = $true

# Column A needs to widen to fit the new, longer row labels.
$ws.Columns.Item(1).ColumnWidth = 23
